# The "optimization_parameters" sheet had a stray leftover row (a row whose
# first cell literally read "Sheet", with two orphaned numbers next to it)
# sitting above the real "simulation_timepoints" row. Clean it up by
# deleting that entire row, which shifts the real data up into its place.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

$ws.Activate()
$ws.Rows.Item(16).EntireRow.Select()
$ws.Rows.Item(16).EntireRow.Delete()

# Finish up on a different sheet, as reflected by the saved workbook state.
$wb.Worksheets.Item("dcin5_log2_optimized_expression").Activate()
